$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custom Table Entry")

$ws.Range("A4").Value = "mesu5p700d3ja57"
$ws.Range("A5").Value = "67rgxmn16wdl1t3"
$ws.Range("A6").Value = "978vx4h4n1d0t12"
$ws.Range("A7").Value = "qtp57advs3xssvx"

$ws.Range("A9").Value = "x00hb59d369oi73"
$ws.Range("A10").Value = "ty2if8tqb09qg33"
$ws.Range("A11").Value = "j91su458oloof53"

$ws.Range("A13").Value = "7e4m3n9on40nvc6"
